$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 8236.454
$ws.Range("I70").Value = 1671.2858
$ws.Range("J70").Value = 19725.5
$ws.Range("K70").Value = 5013.857400000001
$ws.Range("L70").Value = 59176.5
$ws.Range("M70").Value = -4743.857400000001
$ws.Range("N70").Value = -59716.5
$ws.Range("H73").Value = 8236.454
$ws.Range("I73").Value = 1671.2858
$ws.Range("J73").Value = 19725.5
$ws.Range("K73").Value = 5013.857400000001
$ws.Range("L73").Value = 59176.5
$ws.Range("M73").Value = -4077.857400000001
$ws.Range("N73").Value = -61048.5
$ws.Range("H74").Value = 4372.5
$ws.Range("I74").Value = 4166.6665
$ws.Range("K74").Value = 4166.6665
$ws.Range("M74").Value = -3230.6665
$ws.Range("H77").Value = 4372.5
$ws.Range("I77").Value = 4166.6665
$ws.Range("K77").Value = 20833.3325
$ws.Range("M77").Value = -16153.3325
$ws.Range("H98").Value = 1477.1177
$ws.Range("I98").Value = 1581.5385
$ws.Range("K98").Value = 1581.5385
$ws.Range("M98").Value = -83.53850000000011
$ws.Range("H112").Value = 2175.7
$ws.Range("I112").Value = 1075
$ws.Range("K112").Value = 3225
$ws.Range("M112").Value = -2117
$ws.Range("H122").Value = 1477.1177
$ws.Range("I122").Value = 1581.5385
$ws.Range("K122").Value = 4744.6155
$ws.Range("M122").Value = -2294.6155
$ws.Range("H137").Value = 1699.3077
$ws.Range("I137").Value = 1511.625
$ws.Range("J137").Value = 1999.6
$ws.Range("K137").Value = 4534.875
$ws.Range("L137").Value = 5998.799999999999
$ws.Range("M137").Value = -1984.875
$ws.Range("N137").Value = -11098.8
$ws.Range("H138").Value = 3717.6667
$ws.Range("I138").Value = 3030.6775
$ws.Range("J138").Value = 5653.727
$ws.Range("K138").Value = 9092.032499999999
$ws.Range("L138").Value = 16961.181
$ws.Range("M138").Value = -3952.032499999999
$ws.Range("N138").Value = -27241.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 500
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -288
$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1447.8182
$ws.Range("I20").Value = 1347.1666
$ws.Range("J20").Value = 1568.6
$ws.Range("K20").Value = 1347.1666
$ws.Range("L20").Value = 1568.6
$ws.Range("M20").Value = -1100.1666
$ws.Range("N20").Value = -2062.6
$ws.Range("H86").Value = 1743.7778
$ws.Range("I86").Value = 1657.3334
$ws.Range("K86").Value = 1657.3334
$ws.Range("M86").Value = -534.3334
$ws.Range("H89").Value = 1743.7778
$ws.Range("I89").Value = 1657.3334
$ws.Range("K89").Value = 8286.666999999999
$ws.Range("M89").Value = -2670.666999999999
$ws.Range("H99").Value = 2387.9583
$ws.Range("J99").Value = 2587.7856
$ws.Range("L99").Value = 2587.7856
$ws.Range("N99").Value = -5583.7856
$ws.Range("H134").Value = 2801.353
$ws.Range("I134").Value = 1982.3
$ws.Range("K134").Value = 5946.9
$ws.Range("M134").Value = -3411.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 532.25
$ws.Range("I2").Value = 676.3333
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 676.3333
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -563.3333
$ws.Range("N2").Value = -326
$ws.Range("H12").Value = 8333.333000000001
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -9830
$ws.Range("N12").Value = -5340
$ws.Range("H28").Value = 11888.667
$ws.Range("J28").Value = 12833
$ws.Range("L28").Value = 12833
$ws.Range("N28").Value = -13323
$ws.Range("H31").Value = 3336.1304
$ws.Range("I31").Value = 2495.375
$ws.Range("J31").Value = 5257.857
$ws.Range("K31").Value = 2495.375
$ws.Range("L31").Value = 5257.857
$ws.Range("M31").Value = -2200.375
$ws.Range("N31").Value = -5847.857
$ws.Range("H34").Value = 3336.1304
$ws.Range("I34").Value = 2495.375
$ws.Range("J34").Value = 5257.857
$ws.Range("K34").Value = 2495.375
$ws.Range("L34").Value = 5257.857
$ws.Range("M34").Value = -2293.375
$ws.Range("N34").Value = -5661.857
$ws.Range("H58").Value = 2002.3611
$ws.Range("I58").Value = 1141.6538
$ws.Range("J58").Value = 4240.2
$ws.Range("K58").Value = 1141.6538
$ws.Range("L58").Value = 4240.2
$ws.Range("M58").Value = -938.6538
$ws.Range("N58").Value = -4646.2
$ws.Range("H62").Value = 46826.223
$ws.Range("I62").Value = 2548.1428
$ws.Range("J62").Value = 201799.5
$ws.Range("K62").Value = 2548.1428
$ws.Range("L62").Value = 201799.5
$ws.Range("M62").Value = -1924.1428
$ws.Range("N62").Value = -203047.5
$ws.Range("H65").Value = 46826.223
$ws.Range("I65").Value = 2548.1428
$ws.Range("J65").Value = 201799.5
$ws.Range("K65").Value = 12740.714
$ws.Range("L65").Value = 1008997.5
$ws.Range("M65").Value = -9620.714
$ws.Range("N65").Value = -1015237.5
$ws.Range("H68").Value = 44000
$ws.Range("J68").Value = 44000
$ws.Range("L68").Value = 44000
$ws.Range("N68").Value = -45498
$ws.Range("H71").Value = 44000
$ws.Range("J71").Value = 44000
$ws.Range("L71").Value = 132000
$ws.Range("N71").Value = -139488
$ws.Range("H134").Value = 2303.818
$ws.Range("I134").Value = 1983.5555
$ws.Range("K134").Value = 5950.666499999999
$ws.Range("M134").Value = -3415.666499999999
$ws.Range("H136").Value = 2002.3611
$ws.Range("I136").Value = 1141.6538
$ws.Range("J136").Value = 4240.2
$ws.Range("K136").Value = 3424.9614
$ws.Range("L136").Value = 12720.6
$ws.Range("M136").Value = -874.9614000000001
$ws.Range("N136").Value = -17820.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 286.17648
$ws.Range("I92").Value = 276.36365
$ws.Range("J92").Value = 304.16666
$ws.Range("K92").Value = 829.09095
$ws.Range("L92").Value = 912.4999799999999
$ws.Range("M92").Value = 418.90905
$ws.Range("N92").Value = -3408.49998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1633
$ws.Range("I9").Value = 1633
$ws.Range("K9").Value = 1633
$ws.Range("M9").Value = -1463
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338
$ws.Range("H38").Value = 12500
$ws.Range("J38").Value = 12500
$ws.Range("L38").Value = 12500
$ws.Range("N38").Value = -13426
$ws.Range("H70").Value = 6623.5356
$ws.Range("I70").Value = 6157.3184
$ws.Range("K70").Value = 6157.3184
$ws.Range("M70").Value = -5887.3184
$ws.Range("H73").Value = 6623.5356
$ws.Range("I73").Value = 6157.3184
$ws.Range("K73").Value = 6157.3184
$ws.Range("M73").Value = -5221.3184
$ws.Range("H80").Value = 5951.857
$ws.Range("J80").Value = 5951.857
$ws.Range("L80").Value = 5951.857
$ws.Range("N80").Value = -7947.857
$ws.Range("H83").Value = 5951.857
$ws.Range("J83").Value = 5951.857
$ws.Range("L83").Value = 29759.285
$ws.Range("N83").Value = -39743.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2312.8235
$ws.Range("I122").Value = 2664.875
$ws.Range("J122").Value = 1999.8889
$ws.Range("K122").Value = 7994.625
$ws.Range("L122").Value = 5999.6667
$ws.Range("M122").Value = -5544.625
$ws.Range("N122").Value = -10899.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 55135
$ws.Range("J68").Value = 55135
$ws.Range("L68").Value = 55135
$ws.Range("N68").Value = -56757
$ws.Range("H71").Value = 55135
$ws.Range("J71").Value = 55135
$ws.Range("L71").Value = 165405
$ws.Range("N71").Value = -173517
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H126").Value = 3688.2222
$ws.Range("I126").Value = 3199.4
$ws.Range("K126").Value = 9598.200000000001
$ws.Range("M126").Value = -7128.200000000001
$ws.Range("H132").Value = 4009.8
$ws.Range("I132").Value = 3291.3076
$ws.Range("K132").Value = 9873.9228
$ws.Range("M132").Value = -7343.9228
$ws.Range("H136").Value = 1916.7142
$ws.Range("I136").Value = 1283.4
$ws.Range("K136").Value = 3850.2
$ws.Range("M136").Value = -1300.2
